# Regenerate save_data "K" column (was Strike#) with freshly calculated s_vals.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds "K" for rows 2-21 (row 1 is the header).
# Updated / recalculated strikeout (K) values per row.
$kVals = @{
    2  = 1
    3  = 0
    4  = 1
    5  = 0
    6  = 0
    7  = 2
    8  = 0
    9  = 1
    10 = 0
    11 = 1
    12 = 0
    13 = 1
    14 = 2
    15 = 0
    16 = 0
    18 = 2
    19 = 1
    20 = 2
}

foreach ($row in $kVals.Keys) {
    $ws.Cells.Item($row, 7).Value = $kVals[$row]
}
